$d = $word.ActiveDocument

# Splits the Range $range (whose formatted run(s) span exactly $range) into
# separate <w:r> runs by deleting its text and re-inserting the replacement
# OOXML fragment $xmlFragment via InsertXML. This produces clean, separate
# <w:r> elements instead of the single merged run Word normally keeps.
function Split-RunXml($range, $xmlFragment) {
    $start = $range.Start
    $range.Delete()
    $insPoint = $d.Range($start, $start)
    $insPoint.InsertXML($xmlFragment)
}

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- First target: the "Titre1" heading paragraph whose text is "{m:v.name}" ---
# Today its leading "{m" is a single run. Split it into two runs "{" and "m"
# (plain runs, no rPr), leaving the existing ":v.name}" run untouched.
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("{m:v.name}") -and $p.Style.NameLocal -eq "heading 1") {
        $target1 = $p
        break
    }
}
if ($null -eq $target1) {
    throw "Could not find the Titre1 '{m:v.name}' paragraph"
}
$run1 = $d.Range($target1.Range.Start, $target1.Range.Start + 2)
if ($run1.Text -ne "{m") {
    throw "Unexpected text for first target run: [$($run1.Text)]"
}
$xml1 = $pkgHeader + '<w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r></w:p>' + $pkgFooter
Split-RunXml $run1 $xml1

# --- Second target: the paragraph whose text is "{m:endfor}" (lang en-US) ---
# Today its leading "{m:" is a single run. Split it into two runs "{" and
# "m:", both keeping the original <w:rPr><w:lang w:val="en-US"/></w:rPr>.
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("{m:endfor}")) {
        $target2 = $p
        break
    }
}
if ($null -eq $target2) {
    throw "Could not find the '{m:endfor}' paragraph"
}
$run2 = $d.Range($target2.Range.Start, $target2.Range.Start + 3)
if ($run2.Text -ne "{m:") {
    throw "Unexpected text for second target run: [$($run2.Text)]"
}
$xml2 = $pkgHeader + '<w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>m:</w:t></w:r></w:p>' + $pkgFooter
Split-RunXml $run2 $xml2

Write-Output "Split 2 runs successfully"
